$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.269.23"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").Value = "1.908.37"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.725"
$ws.Range("E5").Value = "  +9.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "255.92"
$ws.Range("E6").Value = "  +4.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.72"
$ws.Range("E8").Value = "  -1.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  +8.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.87"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0760"
$ws.Range("E11").Value = "  +5.54%  "

# Row 12
$ws.Range("E12").Value = "  -0.58%  "

# Row 13
$ws.Range("D13").Value = "2.186.32"
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.94"
$ws.Range("E14").Value = "  +7.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.729"
$ws.Range("E15").Value = "  +4.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.98"
$ws.Range("E16").Value = "  +2.91%  "

# Row 17
$ws.Range("D17").Value = "1.903.44"
$ws.Range("E17").Value = "  -0.13%  "

# Row 18
$ws.Range("D18").Value = "35.254.75"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.89"
$ws.Range("E19").Value = "  +3.86%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0852"
$ws.Range("E20").Value = "  +3.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.08"
$ws.Range("E21").Value = "  +1.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.07"
$ws.Range("E22").Value = "  +4.63%  "

# Row 23
$ws.Range("E23").Value = "  +6.06%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +7.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  +3.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.22"
$ws.Range("E27").Value = "  -2.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.72"
$ws.Range("E28").Value = "  +3.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.77"
$ws.Range("E29").Value = "  +2.31%  "

# Row 30
$ws.Range("E30").Value = "  +4.50%  "

# Row 31
$ws.Range("D31").Value = "4.129.32"
$ws.Range("E31").Value = "  +19.47%  "

# Row 32
$ws.Range("E32").Value = "  +6.36%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0591"
$ws.Range("E34").Value = "  +4.66%  "

# Row 35
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +22.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.25"
$ws.Range("E36").Value = "  +3.95%  "

# Row 37
$ws.Range("E37").Value = "  -1.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.914"
$ws.Range("E38").Value = "  -1.93%  "

# Row 39
$ws.Range("E39").Value = "  +0.54%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0219"
$ws.Range("E40").Value = "  +5.67%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.24"
$ws.Range("E41").Value = "  +6.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.51"
$ws.Range("E42").Value = "  +7.77%  "

# Row 43
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0652"
$ws.Range("E44").Value = "  +4.64%  "

# Row 45
$ws.Range("D45").Value = "1.338.24"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  +2.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +0.67%  "

# Row 48
$ws.Range("E48").Value = "  +3.02%  "

# Row 49
$ws.Range("E49").Value = "  -0.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.26"
$ws.Range("E50").Value = "  -5.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0752"
$ws.Range("E51").Value = "  +6.69%  "
